$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "95.736.44"
$ws.Range("E2").Value = "  -0.54%  "
$ws.Range("D3").Value = "3.566.45"
$ws.Range("E3").Value = "  -1.82%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'237.53"
$ws.Range("E5").Value = "  -2.14%  "
$ws.Range("D6").Value = "'653.29"
$ws.Range("E6").Value = "  +0.23%  "
$ws.Range("E7").Value = "  +4.95%  "
$ws.Range("D8").Value = "'0.403"
$ws.Range("E8").Value = "  -0.77%  "
$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").Value = "'1.05"
$ws.Range("E10").Value = "  +2.77%  "
$ws.Range("D11").Value = "3.563.40"
$ws.Range("E11").Value = "  -1.86%  "
$ws.Range("B12").Value = "Avalanche"
$ws.Range("C12").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D12").Value = "'43.00"
$ws.Range("E12").Value = "  -2.00%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'0.202"
$ws.Range("E13").Value = "  +0.39%  "
$ws.Range("D14").Value = "'6.39"
$ws.Range("E14").Value = "  -0.26%  "
$ws.Range("D15").Value = "4.226.33"
$ws.Range("E15").Value = "  -2.13%  "
$ws.Range("D16").Value = "95.475.30"
$ws.Range("E16").Value = "  -0.64%  "
$ws.Range("E17").Value = "  -1.49%  "
$ws.Range("D18").Value = "3.562.33"
$ws.Range("E18").Value = "  -2.05%  "
$ws.Range("B19").Value = "Polkadot"
$ws.Range("C19").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D19").Value = "'8.06"
$ws.Range("E19").Value = "  +3.98%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "'12.59"
$ws.Range("E20").Value = "  -1.61%  "
$ws.Range("D21").Value = "'17.87"
$ws.Range("E21").Value = "  -3.33%  "
$ws.Range("D22").Value = "'0.492"
$ws.Range("E22").Value = "  +2.31%  "
$ws.Range("D23").Value = "'3.44"
$ws.Range("E23").Value = "  -1.07%  "
$ws.Range("D24").Value = "'509.89"
$ws.Range("E24").Value = "  -1.36%  "
$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").Value = "'0.0000200"
$ws.Range("E25").Value = "  +0.94%  "
$ws.Range("B26").Value = "NEARProtocol"
$ws.Range("C26").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D26").Value = "'6.90"
$ws.Range("E26").Value = "  +0.74%  "
$ws.Range("D27").Value = "'92.21"
$ws.Range("E27").Value = "  -5.58%  "
$ws.Range("D28").Value = "'12.72"
$ws.Range("E28").Value = "  -0.63%  "
$ws.Range("D29").Value = "3.752.52"
$ws.Range("E29").Value = "  -1.85%  "
$ws.Range("D30").Value = "'3.02"
$ws.Range("E30").Value = "  -5.43%  "
$ws.Range("E31").Value = "  +2.83%  "
$ws.Range("D32").Value = "'11.58"
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("D33").Value = "'0.997"
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("D34").Value = "'0.997"
$ws.Range("E34").Value = "  +0.19%  "
$ws.Range("D35").Value = "'0.177"
$ws.Range("E35").Value = "  -0.31%  "
$ws.Range("D36").Value = "'31.50"
$ws.Range("E36").Value = "  -2.20%  "
$ws.Range("B37").Value = "PolygonEcosystemToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D37").Value = "'0.567"
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").Value = "'1.61"
$ws.Range("E38").Value = "  +7.92%  "
$ws.Range("D39").Value = "'8.46"
$ws.Range("E39").Value = "  +4.95%  "
$ws.Range("D40").Value = "'590.59"
$ws.Range("E40").Value = "  +4.50%  "
$ws.Range("B41").Value = "USDe"
$ws.Range("C41").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "'0.151"
$ws.Range("E42").Value = "  -1.41%  "
$ws.Range("D43").Value = "'1.86"
$ws.Range("E43").Value = "  +4.26%  "
$ws.Range("D44").Value = "'0.903"
$ws.Range("E44").Value = "  -5.06%  "
$ws.Range("D45").Value = "'5.78"
$ws.Range("E45").Value = "  -1.23%  "
$ws.Range("D46").Value = "'2.29"
$ws.Range("E46").Value = "  -2.09%  "
$ws.Range("B47").Value = "WhiteBITCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D47").Value = "'23.45"
$ws.Range("E47").Value = "  -1.38%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'33.88"
$ws.Range("E48").Value = "  -0.85%  "
$ws.Range("D49").Value = "'0.0416"
$ws.Range("E49").Value = "  -2.69%  "
$ws.Range("D50").Value = "'3.55"
$ws.Range("E50").Value = "  +1.51%  "
$ws.Range("D51").Value = "'8.30"
$ws.Range("E51").Value = "  +0.04%  "
